$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Baja California Mexico": row 4 (Commit/Forecast) — Jun/Q2 cleared,
# Jul-Dec/Q3/Q4/FY turnover values updated (monthly rate dropped from
# 0.133341666666667 to 0.111116666666667, quarters/FY recomputed).
# ---------------------------------------------------------------------------
$wsBaja = $wb.Worksheets.Item("Baja California Mexico")
$wsBaja.Range("M4").ClearContents()
$wsBaja.Range("N4").ClearContents()
$wsBaja.Range("O4").Value = 0.111116666666667
$wsBaja.Range("P4").Value = 0.111116666666667
$wsBaja.Range("Q4").Value = 0.111116666666667
$wsBaja.Range("R4").Value = 0.33335
$wsBaja.Range("S4").Value = 0.111116666666667
$wsBaja.Range("T4").Value = 0.111116666666667
$wsBaja.Range("U4").Value = 0.111116666666667
$wsBaja.Range("V4").Value = 0.33335
$wsBaja.Range("W4").Value = 1.3334

# ---------------------------------------------------------------------------
# Sheet "Charlotte  North Carolina": row 4 (Commit/Forecast) — Jun/Q2 cleared.
# ---------------------------------------------------------------------------
$wsCharlotte = $wb.Worksheets.Item("Charlotte  North Carolina")
$wsCharlotte.Range("M4").ClearContents()
$wsCharlotte.Range("N4").ClearContents()

# ---------------------------------------------------------------------------
# Sheet "Cleveland Ohio": row 6 (Commit/Forecast) — Jun/Q2 cleared.
# ---------------------------------------------------------------------------
$wsCleveland = $wb.Worksheets.Item("Cleveland Ohio")
$wsCleveland.Range("M6").ClearContents()
$wsCleveland.Range("N6").ClearContents()

# ---------------------------------------------------------------------------
# Sheet "Marengo Illinois": row 4 (Commit/Forecast) — Jun/Q2 cleared.
# ---------------------------------------------------------------------------
$wsMarengo = $wb.Worksheets.Item("Marengo Illinois")
$wsMarengo.Range("M4").ClearContents()
$wsMarengo.Range("N4").ClearContents()

# ---------------------------------------------------------------------------
# Sheet "Apodaca Pmc Plant 2 Mexico": ytd (E) cleared on PY Actual / AOP rows,
# and the Commit/Forecast row (row 4) is removed entirely.
# ---------------------------------------------------------------------------
$wsApodaca = $wb.Worksheets.Item("Apodaca Pmc Plant 2 Mexico")
$wsApodaca.Range("E2").ClearContents()
$wsApodaca.Range("E3").ClearContents()
$wsApodaca.Rows.Item(4).Delete()
